$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- v1 - FOR (left block: B/C sorting-machine time, E/F output time) ----
# sortRows
$ws.Range("B14").Value = 0.35779300000000003
$ws.Range("C14").Value = 0.15864500000000001
$ws.Range("E14").Value = 0.47426600000000002
$ws.Range("F14").Value = 0.44431500000000002

# outputSortedRows
$ws.Range("B15").Value = 2.09293
$ws.Range("C15").Value = 0.26318999999999998
$ws.Range("E15").Value = 2.1465299999999998
$ws.Range("F15").Value = 0.31968600000000003

# sortAll
$ws.Range("B16").Value = 0.247889
$ws.Range("C16").Value = 0.16555800000000001
$ws.Range("E16").Value = 0.27794099999999999
$ws.Range("F16").Value = 0.17538999999999999

# outputSortedAll
$ws.Range("B17").Value = 1.74708
$ws.Range("C17").Value = 0.31976199999999999
$ws.Range("E17").Value = 2.4739499999999999
$ws.Range("F17").Value = 0.31193399999999999

# calcMovingAve
$ws.Range("B18").Value = 0.42764400000000002
$ws.Range("C18").Value = 0.074173600000000006
$ws.Range("E18").Value = 0.48702299999999998
$ws.Range("F18").Value = 0.0977272

# outputAveRows
$ws.Range("B19").Value = 1.9334
$ws.Range("C19").Value = 0.276758
$ws.Range("E19").Value = 2.1504099999999999
$ws.Range("F19").Value = 0.338003

# ---- v2 - SECTIONS (right block: K/L sorting-machine time, N/O output time) ----
# sortRows
$ws.Range("K14").Value = 0.095850699999999997
$ws.Range("L14").Value = 0.047730300000000003
$ws.Range("N14").Value = 0.17314599999999999
$ws.Range("O14").Value = 0.102475

# outputSortedRows
$ws.Range("K15").Value = 0.54662500000000003
$ws.Range("L15").Value = 0.098079600000000003
$ws.Range("N15").Value = 1.47139
$ws.Range("O15").Value = 0.16388900000000001

# sortAll
$ws.Range("K16").Value = 0.079160499999999995
$ws.Range("L16").Value = 0.068063899999999997
$ws.Range("N16").Value = 0.077828999999999995
$ws.Range("O16").Value = 0.056844899999999997

# outputSortedAll
$ws.Range("K17").Value = 0.58961699999999995
$ws.Range("L17").Value = 0.10283
$ws.Range("N17").Value = 1.28111
$ws.Range("O17").Value = 0.200625

# calcMovingAve
$ws.Range("K18").Value = 0.112903
$ws.Range("L18").Value = 0.022629799999999999
$ws.Range("N18").Value = 0.25892500000000002
$ws.Range("O18").Value = 0.0434229

# outputAveRows
$ws.Range("K19").Value = 0.62946500000000005
$ws.Range("L19").Value = 0.115315
$ws.Range("N19").Value = 1.2981199999999999
$ws.Range("O19").Value = 0.18610099999999999

# Move the active selection on Sheet1 from K3 to Q21
$ws.Range("Q21").Select()
